$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.068.31"
$ws.Range("E2").Value = "  -7.78%  "
$ws.Range("D3").Value = "1.427.67"
$ws.Range("E3").Value = "  -7.45%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.94"
$ws.Range("E6").Value = "  -5.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3714"
$ws.Range("E7").Value = "  -4.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3079"
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.16"
$ws.Range("E9").Value = "  -7.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.004"
$ws.Range("E10").Value = "  -6.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06585"
$ws.Range("E11").Value = "  -8.56%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  -3.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.23"
$ws.Range("E14").Value = "  -7.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.179"
$ws.Range("E15").Value = "  -6.62%  "
$ws.Range("D16").Value = "1.433.86"
$ws.Range("E16").Value = "  -7.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001012"
$ws.Range("E17").Value = "  -8.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05819"
$ws.Range("E18").Value = "  -11.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.49"
$ws.Range("E19").Value = "  -9.18%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.695"
$ws.Range("E21").Value = "  -7.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.52"
$ws.Range("E22").Value = "  -5.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.11"
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.334"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").Value = "20.084.14"
$ws.Range("E25").Value = "  -7.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.289"
$ws.Range("E26").Value = "  -4.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "138.63"
$ws.Range("E27").Value = "  -3.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.92"
$ws.Range("D29").Value = "1.595.08"
$ws.Range("E29").Value = "  -7.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.68"
$ws.Range("E30").Value = "  -6.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.901"
$ws.Range("E31").Value = "  -19.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9173"
$ws.Range("E32").Value = "  -5.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.440"
$ws.Range("E33").Value = "  -7.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07788"
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.401"
$ws.Range("E35").Value = "  -6.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.48"
$ws.Range("E36").Value = "  +8.18%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.774"
$ws.Range("E38").Value = "  -7.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05689"
$ws.Range("E39").Value = "  -6.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1921"
$ws.Range("E40").Value = "  -6.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.123"
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02017"
$ws.Range("E42").Value = "  -8.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.306"
$ws.Range("E43").Value = "  -9.12%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5336"
$ws.Range("E44").Value = "  -7.69%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.557"
$ws.Range("E45").Value = "  -4.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.21"
$ws.Range("E46").Value = "  -6.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5147"
$ws.Range("E47").Value = "  -7.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.782"
$ws.Range("E48").Value = "  -5.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.66"
$ws.Range("E49").Value = "  -5.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.055"
$ws.Range("E50").Value = "  -6.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  +0.21%  "

Write-Output "done"